# Re-process the data with the newly curated dimensions.
# Updates the "provincia" column metadata block (column D) so that it
# mirrors the structure already used by the other columns:
#   D2: sdmx-dimension:refArea -> iaest-measure:provincia
#   D3: dim                    -> medida
#   D4: URI-Provincia          -> xsd:int

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "iaest-measure:provincia"
$ws.Range("D3").Value = "medida"
$ws.Range("D4").Value = "xsd:int"
